$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.512.34'
$ws.Range('E2').Value = '  +3.10%  '
$ws.Range('D3').Value = '1.840.76'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.49%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.45'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.07%  '
$ws.Range('E9').Value = '  +8.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0700'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.46%  '
$ws.Range('E11').Value = '  +2.91%  '
$ws.Range('D12').Value = '2.107.47'
$ws.Range('E12').Value = '  +2.25%  '
$ws.Range('D13').Value = '1.848.20'
$ws.Range('E13').Value = '  +2.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.30'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.674'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.61%  '
$ws.Range('D17').Value = '35.515.69'
$ws.Range('E17').Value = '  +3.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.67%  '
$ws.Range('D19').Value = '0.0₃0802'
$ws.Range('E19').Value = '  +5.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '244.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.42%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('E24').Value = '  +3.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.122'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('E29').Value = '  +26.81%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').Value = '3.304.71'
$ws.Range('E31').Value = '  +36.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0550'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.77%  '
$ws.Range('E33').Value = '  +7.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.54%  '
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '95.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.692'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.50%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.348.03'
$ws.Range('E39').Value = '  +3.77%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.81%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.45'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.19%  '
$ws.Range('E42').Value = '  +5.79%  '
$ws.Range('E43').Value = '  +7.84%  '
$ws.Range('E44').Value = '  +4.20%  '
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0518'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = '2.012.95'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.31%  '
